# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly scraped counts (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 16358
$wsExhibit.Range("F3").Value = 353
$wsExhibit.Range("F4").Value = 734
$wsExhibit.Range("F5").Value = 253
$wsExhibit.Range("F6").Value = 689
$wsExhibit.Range("F7").Value = 1719
$wsExhibit.Range("F8").Value = 160

# --- Sheet 4: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 16358
$wsAll.Range("F3").Value = 353
$wsAll.Range("F4").Value = 734
$wsAll.Range("F5").Value = 253
$wsAll.Range("F8").Value = 689
$wsAll.Range("F9").Value = 1719
$wsAll.Range("F11").Value = 160
